# Updated cryptos list - applies scraped Price (D) and Volume(1h) (E) changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.409.93"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "1.836.03"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("E4").Value = "  +1.06%  "
$ws.Range("D5").Value = "314.39"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D7").Value = "0.4736"
$ws.Range("E7").Value = "  +2.04%  "
$ws.Range("D8").Value = "0.3696"
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").Value = "0.07461"
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("D10").Value = "'0.8860"
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("D11").Value = "20.47"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "1.882.75"
$ws.Range("E12").Value = "  +3.01%  "
$ws.Range("D13").Value = "0.07345"
$ws.Range("E13").Value = "  +3.47%  "
$ws.Range("D14").Value = "5.458"
$ws.Range("E14").Value = "  +1.68%  "
$ws.Range("D15").Value = "93.14"
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("D16").Value = "6.587"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "'1.010"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "0.000008827"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D20").Value = "27.764.56"
$ws.Range("E20").Value = "  +3.09%  "
$ws.Range("D21").Value = "14.81"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").Value = "5.319"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").Value = "10.69"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").Value = "2.110.52"
$ws.Range("E24").Value = "  +2.58%  "
$ws.Range("D25").Value = "1.895"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "152.31"
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("D28").Value = "2.144"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").Value = "5.249"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").Value = "117.65"
$ws.Range("E30").Value = "  +1.96%  "
$ws.Range("D31").Value = "0.08998"
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("D32").Value = "0.7581"
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("D34").Value = "'4.550"
$ws.Range("E34").Value = "  +1.60%  "
$ws.Range("D35").Value = "2.947"
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("D37").Value = "1.105"
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("D38").Value = "0.05351"
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").Value = "2.987"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").Value = "7.337"
$ws.Range("E41").Value = "  +1.20%  "
$ws.Range("D42").Value = "2.412"
$ws.Range("E42").Value = "  +3.87%  "
$ws.Range("D43").Value = "0.5334"
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("D44").Value = "'0.1660"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").Value = "8.512"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").Value = "0.4921"
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("D47").Value = "10.56"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("D49").Value = "105.08"
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("D50").Value = "1.678"
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("D51").Value = "0.06301"
$ws.Range("E51").Value = "  +0.09%  "

Write-Output "Updated cryptos list"
